$d = $word.ActiveDocument

# 1. Update the date line
$d.Content.Find.Execute("2024-06-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-06 Thursday", 2) | Out-Null

# 2. Restructure the row that gains/loses a cell (table row index 6, 1-based):
#    insert a fresh row with the five final values, then drop the old row.
$t = $d.Tables.Item(1)
$oldRow = $t.Rows.Item(6)
$newRow = $t.Rows.Add($oldRow)
$newRow.Cells.Item(1).Range.Text = "15+62=77"
$newRow.Cells.Item(2).Range.Text = "0+18=18"
$newRow.Cells.Item(3).Range.Text = "96+3=99"
$newRow.Cells.Item(4).Range.Text = "8+6=14"
$newRow.Cells.Item(5).Range.Text = "94-48=46"
$t.Rows.Item(7).Delete()

# 3. Every other answer cell: find-and-replace the old value with the new one.
$d.Content.Find.Execute("81+12=93", $true, $false, $false, $false, $false, $true, 1, $false, "84-70=14", 2) | Out-Null
$d.Content.Find.Execute("8+14=22", $true, $false, $false, $false, $false, $true, 1, $false, "17+28=45", 2) | Out-Null
$d.Content.Find.Execute("96-95=1", $true, $false, $false, $false, $false, $true, 1, $false, "98-42=56", 2) | Out-Null
$d.Content.Find.Execute("3+93=96", $true, $false, $false, $false, $false, $true, 1, $false, "61-49=12", 2) | Out-Null
$d.Content.Find.Execute("95-45=50", $true, $false, $false, $false, $false, $true, 1, $false, "57+40=97", 2) | Out-Null
$d.Content.Find.Execute("44-28=16", $true, $false, $false, $false, $false, $true, 1, $false, "15-5=10", 2) | Out-Null
$d.Content.Find.Execute("39+23=62", $true, $false, $false, $false, $false, $true, 1, $false, "88-34=54", 2) | Out-Null
$d.Content.Find.Execute("15-12=3", $true, $false, $false, $false, $false, $true, 1, $false, "79-45=34", 2) | Out-Null
$d.Content.Find.Execute("48-13=35", $true, $false, $false, $false, $false, $true, 1, $false, "92-29=63", 2) | Out-Null
$d.Content.Find.Execute("95-73=22", $true, $false, $false, $false, $false, $true, 1, $false, "8+59=67", 2) | Out-Null
$d.Content.Find.Execute("67-54=13", $true, $false, $false, $false, $false, $true, 1, $false, "37+15=52", 2) | Out-Null
$d.Content.Find.Execute("85-63=22", $true, $false, $false, $false, $false, $true, 1, $false, "17-7=10", 2) | Out-Null
$d.Content.Find.Execute("57+36=93", $true, $false, $false, $false, $false, $true, 1, $false, "96-86=10", 2) | Out-Null
$d.Content.Find.Execute("43+6=49", $true, $false, $false, $false, $false, $true, 1, $false, "75+4=79", 2) | Out-Null
$d.Content.Find.Execute("40-20=20", $true, $false, $false, $false, $false, $true, 1, $false, "87+11=98", 2) | Out-Null
$d.Content.Find.Execute("23+68=91", $true, $false, $false, $false, $false, $true, 1, $false, "44-13=31", 2) | Out-Null
$d.Content.Find.Execute("84-5=79", $true, $false, $false, $false, $false, $true, 1, $false, "32+28=60", 2) | Out-Null
$d.Content.Find.Execute("72-37=35", $true, $false, $false, $false, $false, $true, 1, $false, "56-7=49", 2) | Out-Null
$d.Content.Find.Execute("67+1=68", $true, $false, $false, $false, $false, $true, 1, $false, "5+50=55", 2) | Out-Null
$d.Content.Find.Execute("63+29=92", $true, $false, $false, $false, $false, $true, 1, $false, "73-59=14", 2) | Out-Null
$d.Content.Find.Execute("31-14=17", $true, $false, $false, $false, $false, $true, 1, $false, "53-11=42", 2) | Out-Null
$d.Content.Find.Execute("86-7=79", $true, $false, $false, $false, $false, $true, 1, $false, "45-15=30", 2) | Out-Null
$d.Content.Find.Execute("69+26=95", $true, $false, $false, $false, $false, $true, 1, $false, "99-26=73", 2) | Out-Null
$d.Content.Find.Execute("67-66=1", $true, $false, $false, $false, $false, $true, 1, $false, "16+1=17", 2) | Out-Null
$d.Content.Find.Execute("87-71=16", $true, $false, $false, $false, $false, $true, 1, $false, "81+7=88", 2) | Out-Null
$d.Content.Find.Execute("31-15=16", $true, $false, $false, $false, $false, $true, 1, $false, "74-54=20", 2) | Out-Null
$d.Content.Find.Execute("76-25=51", $true, $false, $false, $false, $false, $true, 1, $false, "95-41=54", 2) | Out-Null
$d.Content.Find.Execute("50+46=96", $true, $false, $false, $false, $false, $true, 1, $false, "70-60=10", 2) | Out-Null
$d.Content.Find.Execute("1+58=59", $true, $false, $false, $false, $false, $true, 1, $false, "53+25=78", 2) | Out-Null
$d.Content.Find.Execute("46+45=91", $true, $false, $false, $false, $false, $true, 1, $false, "36+32=68", 2) | Out-Null
$d.Content.Find.Execute("45-26=19", $true, $false, $false, $false, $false, $true, 1, $false, "20+61=81", 2) | Out-Null
$d.Content.Find.Execute("35+13=48", $true, $false, $false, $false, $false, $true, 1, $false, "96-80=16", 2) | Out-Null
$d.Content.Find.Execute("56+4=60", $true, $false, $false, $false, $false, $true, 1, $false, "86+0=86", 2) | Out-Null
$d.Content.Find.Execute("35-29=6", $true, $false, $false, $false, $false, $true, 1, $false, "31+40=71", 2) | Out-Null
$d.Content.Find.Execute("10+76=86", $true, $false, $false, $false, $false, $true, 1, $false, "8+60=68", 2) | Out-Null
$d.Content.Find.Execute("7+22=29", $true, $false, $false, $false, $false, $true, 1, $false, "43+1=44", 2) | Out-Null
$d.Content.Find.Execute("50-15=35", $true, $false, $false, $false, $false, $true, 1, $false, "78-18=60", 2) | Out-Null
$d.Content.Find.Execute("26-23=3", $true, $false, $false, $false, $false, $true, 1, $false, "47-5=42", 2) | Out-Null
$d.Content.Find.Execute("27-21=6", $true, $false, $false, $false, $false, $true, 1, $false, "7+26=33", 2) | Out-Null
$d.Content.Find.Execute("2+94=96", $true, $false, $false, $false, $false, $true, 1, $false, "25+50=75", 2) | Out-Null
$d.Content.Find.Execute("69-57=12", $true, $false, $false, $false, $false, $true, 1, $false, "68-64=4", 2) | Out-Null
$d.Content.Find.Execute("83+10=93", $true, $false, $false, $false, $false, $true, 1, $false, "14-8=6", 2) | Out-Null
$d.Content.Find.Execute("99-57=42", $true, $false, $false, $false, $false, $true, 1, $false, "94-73=21", 2) | Out-Null
$d.Content.Find.Execute("40+18=58", $true, $false, $false, $false, $false, $true, 1, $false, "45+24=69", 2) | Out-Null
$d.Content.Find.Execute("11+63=74", $true, $false, $false, $false, $false, $true, 1, $false, "40+51=91", 2) | Out-Null
$d.Content.Find.Execute("95-88=7", $true, $false, $false, $false, $false, $true, 1, $false, "40+31=71", 2) | Out-Null
$d.Content.Find.Execute("90-65=25", $true, $false, $false, $false, $false, $true, 1, $false, "31+34=65", 2) | Out-Null
$d.Content.Find.Execute("71-51=20", $true, $false, $false, $false, $false, $true, 1, $false, "8+68=76", 2) | Out-Null
$d.Content.Find.Execute("28+34=62", $true, $false, $false, $false, $false, $true, 1, $false, "85+13=98", 2) | Out-Null
$d.Content.Find.Execute("18+64=82", $true, $false, $false, $false, $false, $true, 1, $false, "21+52=73", 2) | Out-Null
$d.Content.Find.Execute("29+21=50", $true, $false, $false, $false, $false, $true, 1, $false, "23+28=51", 2) | Out-Null
$d.Content.Find.Execute("92-27=65", $true, $false, $false, $false, $false, $true, 1, $false, "93-53=40", 2) | Out-Null
$d.Content.Find.Execute("57+34=91", $true, $false, $false, $false, $false, $true, 1, $false, "5+87=92", 2) | Out-Null
$d.Content.Find.Execute("5+13=18", $true, $false, $false, $false, $false, $true, 1, $false, "41-27=14", 2) | Out-Null
$d.Content.Find.Execute("51-2=49", $true, $false, $false, $false, $false, $true, 1, $false, "15+6=21", 2) | Out-Null
$d.Content.Find.Execute("66-35=31", $true, $false, $false, $false, $false, $true, 1, $false, "58+9=67", 2) | Out-Null
$d.Content.Find.Execute("92-7=85", $true, $false, $false, $false, $false, $true, 1, $false, "14+21=35", 2) | Out-Null
$d.Content.Find.Execute("47-26=21", $true, $false, $false, $false, $false, $true, 1, $false, "10-7=3", 2) | Out-Null
$d.Content.Find.Execute("79-48=31", $true, $false, $false, $false, $false, $true, 1, $false, "53-51=2", 2) | Out-Null
$d.Content.Find.Execute("12+0=12", $true, $false, $false, $false, $false, $true, 1, $false, "55-33=22", 2) | Out-Null
$d.Content.Find.Execute("70-26=44", $true, $false, $false, $false, $false, $true, 1, $false, "99-18=81", 2) | Out-Null
$d.Content.Find.Execute("29-0=29", $true, $false, $false, $false, $false, $true, 1, $false, "25-13=12", 2) | Out-Null
$d.Content.Find.Execute("45-33=12", $true, $false, $false, $false, $false, $true, 1, $false, "5+50=55", 2) | Out-Null
$d.Content.Find.Execute("8+9=17", $true, $false, $false, $false, $false, $true, 1, $false, "14+23=37", 2) | Out-Null
$d.Content.Find.Execute("81-45=36", $true, $false, $false, $false, $false, $true, 1, $false, "80-62=18", 2) | Out-Null
$d.Content.Find.Execute("29+51=80", $true, $false, $false, $false, $false, $true, 1, $false, "79-38=41", 2) | Out-Null
$d.Content.Find.Execute("79-75=4", $true, $false, $false, $false, $false, $true, 1, $false, "50-22=28", 2) | Out-Null
$d.Content.Find.Execute("35+51=86", $true, $false, $false, $false, $false, $true, 1, $false, "57-50=7", 2) | Out-Null
$d.Content.Find.Execute("82-3=79", $true, $false, $false, $false, $false, $true, 1, $false, "22-14=8", 2) | Out-Null
$d.Content.Find.Execute("66-32=34", $true, $false, $false, $false, $false, $true, 1, $false, "30+37=67", 2) | Out-Null
$d.Content.Find.Execute("78+16=94", $true, $false, $false, $false, $false, $true, 1, $false, "91-52=39", 2) | Out-Null
$d.Content.Find.Execute("82-68=14", $true, $false, $false, $false, $false, $true, 1, $false, "47+31=78", 2) | Out-Null
$d.Content.Find.Execute("85-78=7", $true, $false, $false, $false, $false, $true, 1, $false, "32+5=37", 2) | Out-Null
$d.Content.Find.Execute("96-68=28", $true, $false, $false, $false, $false, $true, 1, $false, "9+8=17", 2) | Out-Null
$d.Content.Find.Execute("36+63=99", $true, $false, $false, $false, $false, $true, 1, $false, "37+1=38", 2) | Out-Null
$d.Content.Find.Execute("25+34=59", $true, $false, $false, $false, $false, $true, 1, $false, "67-5=62", 2) | Out-Null
$d.Content.Find.Execute("44-12=32", $true, $false, $false, $false, $false, $true, 1, $false, "22+51=73", 2) | Out-Null
$d.Content.Find.Execute("75-48=27", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=32", 2) | Out-Null
$d.Content.Find.Execute("6+38=44", $true, $false, $false, $false, $false, $true, 1, $false, "27-17=10", 2) | Out-Null
$d.Content.Find.Execute("53-43=10", $true, $false, $false, $false, $false, $true, 1, $false, "16+31=47", 2) | Out-Null
$d.Content.Find.Execute("99-75=24", $true, $false, $false, $false, $false, $true, 1, $false, "58-7=51", 2) | Out-Null
$d.Content.Find.Execute("37+49=86", $true, $false, $false, $false, $false, $true, 1, $false, "81-25=56", 2) | Out-Null
$d.Content.Find.Execute("90-53=37", $true, $false, $false, $false, $false, $true, 1, $false, "76-12=64", 2) | Out-Null
$d.Content.Find.Execute("29+42=71", $true, $false, $false, $false, $false, $true, 1, $false, "4+23=27", 2) | Out-Null
$d.Content.Find.Execute("82-19=63", $true, $false, $false, $false, $false, $true, 1, $false, "99-43=56", 2) | Out-Null
$d.Content.Find.Execute("3+54=57", $true, $false, $false, $false, $false, $true, 1, $false, "5+59=64", 2) | Out-Null
$d.Content.Find.Execute("23+7=30", $true, $false, $false, $false, $false, $true, 1, $false, "62+29=91", 2) | Out-Null
$d.Content.Find.Execute("5+7=12", $true, $false, $false, $false, $false, $true, 1, $false, "18-7=11", 2) | Out-Null
$d.Content.Find.Execute("56-24=32", $true, $false, $false, $false, $false, $true, 1, $false, "97+2=99", 2) | Out-Null
$d.Content.Find.Execute("13+4=17", $true, $false, $false, $false, $false, $true, 1, $false, "3+75=78", 2) | Out-Null
$d.Content.Find.Execute("66-58=8", $true, $false, $false, $false, $false, $true, 1, $false, "84-39=45", 2) | Out-Null
$d.Content.Find.Execute("97-70=27", $true, $false, $false, $false, $false, $true, 1, $false, "33+22=55", 2) | Out-Null
$d.Content.Find.Execute("94-16=78", $true, $false, $false, $false, $false, $true, 1, $false, "75+7=82", 2) | Out-Null
$d.Content.Find.Execute("92-87=5", $true, $false, $false, $false, $false, $true, 1, $false, "96-24=72", 2) | Out-Null
$d.Content.Find.Execute("97-83=14", $true, $false, $false, $false, $false, $true, 1, $false, "16+72=88", 2) | Out-Null

Write-Host "done"
